$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.345.37"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.931.83"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.50"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7179"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3266"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("E9").Value = "  +4.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07189"
$ws.Range("E10").Value = "  +5.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7991"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08095"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.930.69"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.423"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.65"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.77"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.324.17"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "251.98"
$ws.Range("E18").Value = "  -3.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008133"
$ws.Range("E19").Value = "  +2.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.792"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.184.48"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.920"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.721"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.24"
$ws.Range("E26").Value = "  +3.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.22"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.323"
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1279"
$ws.Range("E29").Value = "  -3.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.361"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.545"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.426"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.194"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05210"
$ws.Range("E34").Value = "  +2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.267"
$ws.Range("E35").Value = "  +6.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7469"
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.770"
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.803"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.99"
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.441"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4524"
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.028"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8408"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.75"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.782"
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.406"
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.53"
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06074"
$ws.Range("E50").Value = "  +2.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4174"
$ws.Range("E51").Value = "  +1.82%  "

Write-Host "Applied cryptos update"